$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Remove old columns F and G entirely (they no longer exist in the new layout)
$ws.Columns.Item(6).Delete()
$ws.Columns.Item(6).Delete()

# --- Data rows 2-11 ---
$data = @(
    @(0, 7412.066666666667, 8011, 6888, 0.1465380112330119),
    @(1, 7914.2,            8561, 7472, 0.1754698673884074),
    @(2, 8134.533333333334, 8709, 7353, 0.1218748172124227),
    @(3, 8661.866666666667, 9485, 7842, 0.1320118347803752),
    @(4, 7755.366666666667, 8157, 7334, 0.1488465309143066),
    @(5, 7961.333333333333, 8621, 7428, 0.137266739209493),
    @(6, 7914.7,            8429, 7129, 0.1110277016957601),
    @(7, 7682.966666666666, 8354, 6682, 0.1425660371780396),
    @(8, 7146.166666666667, 7599, 6387, 0.1276486317316691),
    @(9, 7523.2,            7979, 6389, 0.1220884879430135)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
